$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the hookup-wire description in the "REQUIRED FOR MID-FEBRUARY" section
$ws.Range("A14").Value = "Hookup wire (stranded is best, any color ok)"

# 2) Insert two new rows for the accelerometer and RGB light sensor, right before
#    the "Any other sensors..." note (old row 32), inside the
#    "PICK AT LEAST TWO OF THE FOLLOWING" section.
$ws.Rows("32:33").Insert()

# Row 32: 3-axis accelerometer
$ws.Range("A32").Value = "3-axis accelerometer"
$ws.Range("B32").Value = 10.49
$ws.Range("C32").Value = 1
$ws.Range("D32").Formula = "=C32*B32"
$ws.Range("E32").Value = "https://www.sparkfun.com/products/13926"

# Row 33: RGB light sensor
$ws.Range("A33").Value = "RGB light sensor"
$ws.Range("B33").Value = 7.95
$ws.Range("C33").Value = 1
$ws.Range("D33").Formula = "=C33*B33"
$ws.Range("E33").Value = "https://www.sparkfun.com/products/12829"

# 3) Restore the selected cell as in the authored workbook
$ws.Range("A15").Select()
